$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Split the paragraph "...learn how to become a doctor right?" so
#    that "right?" becomes its own paragraph (keeping the paragraph
#    formatting that the original single paragraph carried).
# ------------------------------------------------------------------
$needle = "to become a doctor right?"
$findRng = $d.Content
$found = $findRng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'right?' sentence to split."
}
$splitPos = $findRng.End - 6  # length of the word "right?"
$splitRng = $d.Range($splitPos, $splitPos)
$splitRng.InsertParagraphBefore()

# ------------------------------------------------------------------
# 2) The hidden "_GoBack" bookmark used to sit right before the
#    "right?" run; after the split it now opens the new "right?"
#    paragraph. Remove it from there - it is going to be relocated
#    onto the paragraph that used to hold the OLE/ActiveX object.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 3) Find the paragraph that used to contain the embedded OLE/ActiveX
#    control (its Range.Text reads empty because the object itself
#    isn't text) and strip that content out, merging the now-empty
#    paragraph into the final paragraph of the document.
# ------------------------------------------------------------------
$oleParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq [char]13 -and $p.Format.Shading.BackgroundPatternColor -eq 16118770) {
        $oleParaIndex = $i
    }
}
if ($oleParaIndex -eq -1) {
    throw "Could not locate the OLE object paragraph."
}
$oleP = $d.Paragraphs.Item($oleParaIndex)
$mergeRng = $d.Range($oleP.Range.End - 1, $oleP.Range.End)
$mergeRng.Delete()

# ------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark, collapsed, at the start of the
#    (now final) trailing paragraph.
# ------------------------------------------------------------------
$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmRng = $lastP.Range.Duplicate
$bmRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRng)
